$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,5).Value = "2026-02-08 20:48:21"  # E2
$ws.Cells.Item(2,8).Value = "'88%"  # H2
$ws.Cells.Item(2,9).Value = "5.6 mm"  # I2
$ws.Cells.Item(3,5).Value = "2026-02-08 20:48:23"  # E3
$ws.Cells.Item(3,9).Value = "2.2 mm"  # I3
$ws.Cells.Item(4,5).Value = "2026-02-08 20:48:25"  # E4
$ws.Cells.Item(4,8).Value = "'69%"  # H4
$ws.Cells.Item(4,10).Value = "1002.3 hPa"  # J4
$ws.Cells.Item(5,5).Value = "2026-02-08 20:48:28"  # E5
$ws.Cells.Item(6,5).Value = "2026-02-08 20:48:30"  # E6
$ws.Cells.Item(6,10).Value = "1002.2 hPa"  # J6
$ws.Cells.Item(7,5).Value = "2026-02-08 20:48:33"  # E7
$ws.Cells.Item(7,10).Value = "1002.5 hPa"  # J7
$ws.Cells.Item(7,11).Value = "11.2 MJ/m2"  # K7
$ws.Cells.Item(8,5).Value = "2026-02-08 20:48:35"  # E8
$ws.Cells.Item(8,10).Value = "1002.4 hPa"  # J8
$ws.Cells.Item(9,5).Value = "2026-02-08 20:48:37"  # E9
$ws.Cells.Item(10,5).Value = "2026-02-08 20:48:40"  # E10
$ws.Cells.Item(10,9).Value = "3.7 mm"  # I10
$ws.Cells.Item(11,5).Value = "2026-02-08 20:48:42"  # E11
$ws.Cells.Item(12,5).Value = "2026-02-08 20:48:44"  # E12
$ws.Cells.Item(12,8).Value = "'73%"  # H12
$ws.Cells.Item(13,5).Value = "2026-02-08 20:48:47"  # E13
$ws.Cells.Item(13,10).Value = "1003.8 hPa"  # J13
$ws.Cells.Item(14,5).Value = "2026-02-08 20:48:49"  # E14
$ws.Cells.Item(15,5).Value = "2026-02-08 20:48:52"  # E15
$ws.Cells.Item(16,5).Value = "2026-02-08 20:48:54"  # E16
$ws.Cells.Item(16,9).Value = "3.6 mm"  # I16
$ws.Cells.Item(17,5).Value = "2026-02-08 20:48:56"  # E17
$ws.Cells.Item(18,5).Value = "2026-02-08 20:48:59"  # E18
$ws.Cells.Item(18,8).Value = "'70%"  # H18
$ws.Cells.Item(18,10).Value = "1002.5 hPa"  # J18
$ws.Cells.Item(19,5).Value = "2026-02-08 20:49:01"  # E19
$ws.Cells.Item(20,5).Value = "2026-02-08 20:49:07"  # E20
$ws.Cells.Item(20,9).Value = "9.2 mm"  # I20
$ws.Cells.Item(21,5).Value = "2026-02-08 20:49:10"  # E21
$ws.Cells.Item(21,10).Value = "1003.2 hPa"  # J21
$ws.Cells.Item(22,5).Value = "2026-02-08 20:49:12"  # E22
$ws.Cells.Item(23,5).Value = "2026-02-08 20:49:15"  # E23
$ws.Cells.Item(23,9).Value = "5.5 mm"  # I23
$ws.Cells.Item(24,5).Value = "2026-02-08 20:49:17"  # E24
$ws.Cells.Item(24,10).Value = "1003.8 hPa"  # J24
$ws.Cells.Item(25,5).Value = "2026-02-08 20:49:20"  # E25
$ws.Cells.Item(25,9).Value = "0.8 mm"  # I25
$ws.Cells.Item(26,5).Value = "2026-02-08 20:49:22"  # E26
$ws.Cells.Item(26,8).Value = "'70%"  # H26
$ws.Cells.Item(26,10).Value = "1001.6 hPa"  # J26
$ws.Cells.Item(26,15).Value = "3.6 °C"  # O26
$ws.Cells.Item(27,5).Value = "2026-02-08 20:49:25"  # E27
$ws.Cells.Item(28,5).Value = "2026-02-08 20:49:27"  # E28
$ws.Cells.Item(28,8).Value = "'70%"  # H28
$ws.Cells.Item(28,10).Value = "1002.2 hPa"  # J28
$ws.Cells.Item(28,15).Value = "8.5 °C"  # O28
$ws.Cells.Item(29,5).Value = "2026-02-08 20:49:30"  # E29
$ws.Cells.Item(30,5).Value = "2026-02-08 20:49:32"  # E30
$ws.Cells.Item(30,8).Value = "'70%"  # H30
$ws.Cells.Item(30,10).Value = "1002.5 hPa"  # J30
$ws.Cells.Item(31,5).Value = "2026-02-08 20:49:35"  # E31
$ws.Cells.Item(31,8).Value = "'76%"  # H31
$ws.Cells.Item(31,9).Value = "0.7 mm"  # I31
$ws.Cells.Item(31,10).Value = "1001.7 hPa"  # J31
$ws.Cells.Item(31,14).Value = "7.6 °C 20:27 TU"  # N31
$ws.Cells.Item(32,5).Value = "2026-02-08 20:49:37"  # E32
$ws.Cells.Item(33,5).Value = "2026-02-08 20:49:39"  # E33
$ws.Cells.Item(33,10).Value = "1003.4 hPa"  # J33
$ws.Cells.Item(34,5).Value = "2026-02-08 20:49:42"  # E34
$ws.Cells.Item(35,5).Value = "2026-02-08 20:49:44"  # E35
$ws.Cells.Item(35,10).Value = "1004.6 hPa"  # J35
$ws.Cells.Item(36,5).Value = "2026-02-08 20:49:47"  # E36
$ws.Cells.Item(36,10).Value = "1002.6 hPa"  # J36
$ws.Cells.Item(36,15).Value = "11.1 °C"  # O36
$ws.Cells.Item(37,5).Value = "2026-02-08 20:49:49"  # E37
$ws.Cells.Item(37,10).Value = "1003.4 hPa"  # J37
$ws.Cells.Item(38,5).Value = "2026-02-08 20:49:52"  # E38
$ws.Cells.Item(38,9).Value = "4.9 mm"  # I38
$ws.Cells.Item(39,5).Value = "2026-02-08 20:49:54"  # E39
$ws.Cells.Item(39,9).Value = "1.0 mm"  # I39
$ws.Cells.Item(40,5).Value = "2026-02-08 20:49:57"  # E40
$ws.Cells.Item(40,10).Value = "1003.8 hPa"  # J40
$ws.Cells.Item(41,5).Value = "2026-02-08 20:49:59"  # E41
$ws.Cells.Item(41,10).Value = "1002.6 hPa"  # J41
$ws.Cells.Item(42,5).Value = "2026-02-08 20:50:01"  # E42
$ws.Cells.Item(43,5).Value = "2026-02-08 20:50:04"  # E43
$ws.Cells.Item(43,15).Value = "7.1 °C"  # O43
$ws.Cells.Item(44,5).Value = "2026-02-08 20:50:06"  # E44
$ws.Cells.Item(45,5).Value = "2026-02-08 20:50:09"  # E45
$ws.Cells.Item(45,8).Value = "'78%"  # H45
$ws.Cells.Item(45,10).Value = "1004.7 hPa"  # J45
$ws.Cells.Item(46,5).Value = "2026-02-08 20:50:11"  # E46
$ws.Cells.Item(46,10).Value = "1004.3 hPa"  # J46
